$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.866.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.46%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.910.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.18%  "

# Row 4
$ws.Range("E4").Value = "  -0.32%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.43%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3791"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.27%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07279"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.24%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.46%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9046"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.04%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07637"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.43%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.894.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.91%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.475"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.14%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.53%  "

# Row 16
$ws.Range("E16").Value = "  -0.44%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008730"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.66%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.0000"
$ws.Range("D18").Style = "Normal"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "27.894.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.49%  "

# Row 20
$ws.Range("E20").Value = "  -0.53%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.173"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.47%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.123.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.608"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.52%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.09%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.840"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.223"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.82%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.53%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.85%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.895"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.71%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08983"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.48%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.190"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.39%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7897"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.58%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.238"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.27%  "

# Row 35
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.816"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.22%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.668"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.19%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02082"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.21%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.059"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.93%  "

# Row 39
$ws.Range("E39").Value = "  -1.21%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5521"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.40%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05304"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.54%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.787"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.11%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.51%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.490"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1513"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.92%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.67%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4801"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.36%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.22%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.638"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.44%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.23%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06035"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.63%  "
